# Heston benchmark workbook: add a "divr" (dividend rate) parameter column
# to the Param sheet, append three new test cases (rows 9-11) with two
# supporting reference hyperlinks, and add three new worksheets ("8", "9",
# "10") holding the corresponding strike/price (and put-call) tables.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Param")

# ---------------------------------------------------------------------
# 1) Create rows 9, 10, 11 as copies of row 8 (so the new rows inherit the
#    same per-cell number/fill styles already used by the existing data
#    rows), then fill in the actual values for the three new test cases.
# ---------------------------------------------------------------------
$ws.Rows.Item(8).Copy() | Out-Null
$ws.Rows.Item(9).Insert()
$ws.Rows.Item(8).Copy() | Out-Null
$ws.Rows.Item(10).Insert()
$ws.Rows.Item(8).Copy() | Out-Null
$ws.Rows.Item(11).Insert()

# ---------------------------------------------------------------------
# 2) Insert the new "divr" column at J (pushes old df/nc/col_name/Reference
#    from J:M to K:N).
# ---------------------------------------------------------------------
$ws.Columns.Item(10).Insert()

# ---------------------------------------------------------------------
# 3) Header row.
# ---------------------------------------------------------------------
$ws.Range("J1").Value = "divr"

# ---------------------------------------------------------------------
# 4) Fill the "divr" column (new) for the pre-existing rows 2-8 with 0.
# ---------------------------------------------------------------------
for ($r = 2; $r -le 8; $r++) {
    $ws.Cells.Item($r, 10).Value = 0
}

# ---------------------------------------------------------------------
# 5) Row 9 - new test case #8.
# ---------------------------------------------------------------------
$ws.Range("A9").Value = 8
$ws.Range("B9").Value = 0.04
$ws.Range("C9").Value = 0.25
$ws.Range("D9").Value = 1
$ws.Range("E9").Value = -0.5
$ws.Range("F9").Value = 4
$ws.Range("G9").Value = 1
$ws.Range("H9").Value = 100
$ws.Range("I9").Value = 0.01
$ws.Range("J9").Value = 0.02
$ws.Range("K9").Formula = "=4*C9*F9/(D9*D9)"
$ws.Range("L9").Formula = "=2*B9*F9*EXP(-F9*G9)/(D9*D9)/(1-EXP(-F9*G9))"
$ws.Range("M9").Value = "Price"
$ws.Range("N9").Value = "https://financepress.com/2019/02/15/heston-model-reference-prices/"

# ---------------------------------------------------------------------
# 6) Row 10 - new test case #9.
# ---------------------------------------------------------------------
$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 0.01
$ws.Range("C10").Value = 0.25
$ws.Range("D10").Value = 1
$ws.Range("E10").Value = -0.5
$ws.Range("F10").Value = 4
$ws.Range("G10").Value = 0.01
$ws.Range("H10").Value = 100
$ws.Range("I10").Value = 0.01
$ws.Range("J10").Value = 0.02
$ws.Range("K10").Formula = "=4*C10*F10/(D10*D10)"
$ws.Range("L10").Formula = "=2*B10*F10*EXP(-F10*G10)/(D10*D10)/(1-EXP(-F10*G10))"
$ws.Range("M10").Value = "Price"
$ws.Range("N10").Value = "https://financepress.com/2019/02/15/heston-model-reference-prices/"

# ---------------------------------------------------------------------
# 7) Row 11 - new test case #10.
# ---------------------------------------------------------------------
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = 0.114
$ws.Range("C11").Value = 0.043
$ws.Range("D11").Value = 1
$ws.Range("E11").Value = -0.36
$ws.Range("F11").Value = 2.58
$ws.Range("G11").Value = 1
$ws.Range("H11").Value = 100
$ws.Range("I11").Value = 0
$ws.Range("J11").Value = 0
$ws.Range("K11").Formula = "=4*C11*F11/(D11*D11)"
$ws.Range("L11").Formula = "=2*B11*F11*EXP(-F11*G11)/(D11*D11)/(1-EXP(-F11*G11))"
$ws.Range("M11").Value = "Price"
$ws.Range("N11").Value = "Table A3, Figure 3 (Set I) in von Sydow et al. (2018). BENCHOP - SLV. https://doi.org/10.1080/00207160.2018.1544368"
$ws.Range("O11").Value = ""

# ---------------------------------------------------------------------
# 8) Hyperlinks for the financepress reference (rows 9 and 10 only; row 11
#    keeps its reference as plain text).
# ---------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("N9"), "https://financepress.com/2019/02/15/heston-model-reference-prices/") | Out-Null
$ws.Hyperlinks.Add($ws.Range("N10"), "https://financepress.com/2019/02/15/heston-model-reference-prices/") | Out-Null

# ---------------------------------------------------------------------
# 9) Selection / view bookkeeping on Param sheet.
# ---------------------------------------------------------------------
$ws.Range("J11").Select() | Out-Null

# ---------------------------------------------------------------------
# 10) New worksheet "8": Strike / CP / Price table (financepress set I).
# ---------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws9 = $wb.Worksheets.Add($null, $last)
$ws9.Name = "9"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws8 = $wb.Worksheets.Add($null, $last)
$ws8.Name = "8"

$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws10 = $wb.Worksheets.Add($null, $last)
$ws10.Name = "10"

# Reorder so the visible tab order is ... 7, 8, 9, 10
$ws8.Move($ws9, $null)

# --- Sheet "8" data (Strike / CP / Price) ---
$ws8.Range("A1").Value = "Strike"
$ws8.Range("B1").Value = "CP"
$ws8.Range("C1").Value = "Price"

$ws8.Range("A2").Value = 90
$ws8.Range("B2").Value = -1
$ws8.Range("C2").Value = 0.000000045183603586861701

$ws8.Range("A3").Value = 95
$ws8.Range("B3").Value = -1
$ws8.Range("C3").Value = 0.00046195485565385098

$ws8.Range("A4").Value = 100
$ws8.Range("B4").Value = -1
$ws8.Range("C4").Value = 0.47778117162950401

$ws8.Range("A5").Value = 105
$ws8.Range("B5").Value = 1
$ws8.Range("C5").Value = 0.0000025274478231946999

$ws8.Range("A6").Value = 110
$ws8.Range("B6").Value = 1
$ws8.Range("C6").Value = 0.00000000000012993276005262401

$ws8.Range("C2").NumberFormat = "0.00E+00"
$ws8.Range("C5").NumberFormat = "0.00E+00"
$ws8.Range("C6").NumberFormat = "0.00E+00"

$ws8.Range("C6").Select() | Out-Null

# --- Sheet "9" data (Strike / Price) ---
$ws9.Range("A1").Value = "Strike"
$ws9.Range("B1").Value = "Price"

$ws9.Range("A2").Value = 80
$ws9.Range("B2").Value = 7.9588781132567599

$ws9.Range("A3").Value = 90
$ws9.Range("B3").Value = 12.0179667073463

$ws9.Range("A4").Value = 100
$ws9.Range("B4").Value = 17.0552709612701

$ws9.Range("A5").Value = 110
$ws9.Range("B5").Value = 12.1322115167098

$ws9.Range("A6").Value = 120
$ws9.Range("B6").Value = 9.02491348345783

$ws9.Range("C6").Select() | Out-Null

# --- Sheet "10" data (von Sydow et al. moneyness-scaled strikes) ---
$ws10.Range("A1").Value = "Strike"
$ws10.Range("B1").Value = "Price"

$ws10.Range("A2").Formula = "=100*(100/75)"
$ws10.Range("B2").Formula = "=0.908502728459621*(100/75)"

$ws10.Range("A3").Value = 100
$ws10.Range("B3").Formula = "=9.04665011922096"

$ws10.Range("A4").Formula = "=100*(100/125)"
$ws10.Range("B4").Formula = "=28.5147863992987*(100/125)"

$ws10.Range("B2").NumberFormat = "0.000"
$ws10.Range("B3").NumberFormat = "0.000"
$ws10.Range("B4").NumberFormat = "0.000"

$ws10.Range("B10").Select() | Out-Null

# ---------------------------------------------------------------------
# 11) Minor selection bookkeeping on the pre-existing single-table sheets.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("2").Range("B2").Select() | Out-Null
$wb.Worksheets.Item("3").Range("B5").Select() | Out-Null
$wb.Worksheets.Item("4").Range("C15").Select() | Out-Null

$ws.Select()
